# The edit lives in word/styles.xml's <w:docDefaults> block, which isn't
# reachable through the normal Style/Font/ParagraphFormat object model
# (those bind to individual w:style elements, not docDefaults). Word does
# expose the whole package as Flat-OPC XML via Range.WordOpenXML though,
# so round-trip through that to reach the styles part.

$d = $word.ActiveDocument

$full = $d.Content.WordOpenXML

$oldDocDefaults = '<w:docDefaults><w:rPrDefault><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b w:val="0"/><w:i w:val="0"/><w:smallCaps w:val="0"/><w:strike w:val="0"/><w:color w:val="000000"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="none"/><w:shd w:val="clear" w:fill="auto"/><w:vertAlign w:val="baseline"/><w:lang w:val="en"/></w:rPr></w:rPrDefault><w:pPrDefault><w:pPr><w:keepNext w:val="0"/><w:keepLines w:val="0"/><w:widowControl/><w:pBdr><w:top w:val="nil" w:sz="0" w:space="0"/><w:left w:val="nil" w:sz="0" w:space="0"/><w:bottom w:val="nil" w:sz="0" w:space="0"/><w:right w:val="nil" w:sz="0" w:space="0"/><w:between w:val="nil" w:sz="0" w:space="0"/></w:pBdr><w:shd w:val="clear" w:fill="auto"/><w:spacing w:before="0" w:after="0" w:line="276" w:lineRule="auto"/><w:ind w:left="0" w:right="0" w:firstLine="0"/><w:contextualSpacing w:val="0"/><w:jc w:val="left"/></w:pPr></w:pPrDefault></w:docDefaults>'

$newDocDefaults = '<w:docDefaults><w:rPrDefault><w:rPr><w:rFonts w:ascii="Arial" w:cs="Arial" w:eastAsia="Arial" w:hAnsi="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en"/></w:rPr></w:rPrDefault><w:pPrDefault><w:pPr><w:spacing w:line="276" w:lineRule="auto"/></w:pPr></w:pPrDefault></w:docDefaults>'

if ($full.Contains($oldDocDefaults)) {
    $full = $full.Replace($oldDocDefaults, $newDocDefaults)
} else {
    $full = [System.Text.RegularExpressions.Regex]::Replace(
        $full,
        '<w:docDefaults>.*?</w:docDefaults>',
        [System.Text.RegularExpressions.Regex]::Escape($newDocDefaults).Replace('\', '\\'),
        [System.Text.RegularExpressions.RegexOptions]::Singleline
    )
}

$d.Content.WordOpenXML = $full

Write-Output "docDefaults rewritten"
